$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 51
$ws.Range("I11").Value = 51
$ws.Range("K11").Value = 51
$ws.Range("M11").Value = 89

$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = ""

$ws.Range("H46").Value = 995
$ws.Range("J46").Value = 995
$ws.Range("L46").Value = 2985
$ws.Range("N46").Value = -3223

$ws.Range("H47").Value = 2500
$ws.Range("I47").Value = 2500
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 2500
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -1528
$ws.Range("N47").Value = ""

$ws.Range("H60").Value = 995
$ws.Range("J60").Value = 995
$ws.Range("L60").Value = 2985
$ws.Range("N60").Value = -3953

$ws.Range("H80").Value = 1378.6666
$ws.Range("I80").Value = 1120
$ws.Range("J80").Value = 1508
$ws.Range("K80").Value = 3360
$ws.Range("L80").Value = 4524
$ws.Range("M80").Value = -2362
$ws.Range("N80").Value = -6520

$ws.Range("H83").Value = 1378.6666
$ws.Range("I83").Value = 1120
$ws.Range("J83").Value = 1508
$ws.Range("K83").Value = 10080
$ws.Range("L83").Value = 13572
$ws.Range("M83").Value = -5088
$ws.Range("N83").Value = -23556

$ws.Range("H137").Value = 4666.6665
$ws.Range("I137").Value = 4666.6665
$ws.Range("K137").Value = 13999.9995
$ws.Range("M137").Value = -11449.9995

$ws.Range("H138").Value = 3383.182
$ws.Range("J138").Value = 3193.7368
$ws.Range("L138").Value = 9581.2104
$ws.Range("N138").Value = -19861.2104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 100.666664
$ws.Range("I5").Value = 94.333336
$ws.Range("K5").Value = 94.333336
$ws.Range("M5").Value = 17.666664

$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 100.666664
$ws.Range("I4").Value = 94.333336
$ws.Range("K4").Value = 94.333336
$ws.Range("M4").Value = 20.666664

$ws.Range("H22").Value = 234.8
$ws.Range("I22").Value = 234.8
$ws.Range("K22").Value = 234.8
$ws.Range("M22").Value = -61.80000000000001

$ws.Range("H93").Value = 42775
$ws.Range("J93").Value = 42775
$ws.Range("L93").Value = 42775
$ws.Range("N93").Value = -46519

$ws.Range("H134").Value = 2220.1428
$ws.Range("I134").Value = 1479.1666
$ws.Range("K134").Value = 4437.4998
$ws.Range("M134").Value = -1902.4998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 714.6667
$ws.Range("I22").Value = 197.25
$ws.Range("K22").Value = 197.25
$ws.Range("M22").Value = 152.75

$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = ""

$ws.Range("H32").Value = 3178.0908
$ws.Range("I32").Value = 997.75
$ws.Range("K32").Value = 997.75
$ws.Range("M32").Value = -681.75

$ws.Range("H59").Value = 31949.166
$ws.Range("J59").Value = 38425
$ws.Range("L59").Value = 38425
$ws.Range("N59").Value = -40715

$ws.Range("H92").Value = 38000
$ws.Range("J92").Value = 38000
$ws.Range("L92").Value = 38000
$ws.Range("N92").Value = -42992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 200409.8
$ws.Range("I4").Value = 333749.66
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 1001248.98
$ws.Range("L4").Value = 1200
$ws.Range("M4").Value = -1001136.98
$ws.Range("N4").Value = -1424

$ws.Range("H57").Value = 100
$ws.Range("I57").Value = 100
$ws.Range("K57").Value = 300
$ws.Range("M57").Value = 259

$ws.Range("H63").Value = 1750.5
$ws.Range("I63").Value = 1750.5
$ws.Range("K63").Value = 5251.5
$ws.Range("M63").Value = -4502.5

$ws.Range("H64").Value = 912
$ws.Range("I64").Value = 912
$ws.Range("K64").Value = 2736
$ws.Range("M64").Value = -2466

$ws.Range("H66").Value = 1750.5
$ws.Range("I66").Value = 1750.5
$ws.Range("K66").Value = 15754.5
$ws.Range("M66").Value = -12010.5

$ws.Range("H67").Value = 912
$ws.Range("I67").Value = 912
$ws.Range("K67").Value = 2736
$ws.Range("M67").Value = -1800

$ws.Range("H68").Value = 850
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -1289
$ws.Range("N68").Value = -4622

$ws.Range("H71").Value = 850
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 6300
$ws.Range("L71").Value = 9000
$ws.Range("M71").Value = -2244
$ws.Range("N71").Value = -17112

$ws.Range("H103").Value = 3624.5
$ws.Range("J103").Value = 4959.4
$ws.Range("L103").Value = 14878.2
$ws.Range("N103").Value = -16636.2

$ws.Range("H140").Value = 1932
$ws.Range("I140").Value = 1932
$ws.Range("K140").Value = 5796
$ws.Range("M140").Value = -616

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 412.4
$ws.Range("I97").Value = 371.75
$ws.Range("J97").Value = 575
$ws.Range("K97").Value = 371.75
$ws.Range("L97").Value = 575
$ws.Range("M97").Value = 124.25
$ws.Range("N97").Value = -1567

$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

$ws.Range("H132").Value = 2806.2
$ws.Range("I132").Value = 1882.875
$ws.Range("K132").Value = 5648.625
$ws.Range("M132").Value = -3118.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = ""

$ws.Range("H101").Value = 13407
$ws.Range("J101").Value = 13407
$ws.Range("L101").Value = 13407
$ws.Range("N101").Value = -19897

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = ""

$ws.Range("H122").Value = 2535
$ws.Range("I122").Value = 1195
$ws.Range("K122").Value = 3585
$ws.Range("M122").Value = -1135

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -88

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").Value = ""

$ws.Range("H62").Value = 4998.8

$ws.Range("H65").Value = 4998.8

$ws.Range("H82").Value = 39999.5
$ws.Range("J82").Value = 39999.5
$ws.Range("L82").Value = 39999.5
$ws.Range("N82").Value = -40765.5

$ws.Range("H85").Value = 39999.5
$ws.Range("J85").Value = 39999.5
$ws.Range("L85").Value = 39999.5
$ws.Range("N85").Value = -42651.5

$ws.Range("H96").Value = 2999
$ws.Range("I96").Value = 2999
$ws.Range("K96").Value = 2999
$ws.Range("M96").Value = -1626

$ws.Range("H122").Value = 4988.75
$ws.Range("I122").Value = 4800
$ws.Range("K122").Value = 14400
$ws.Range("M122").Value = -11950

